$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Yueling%Ma%NULL%1,               Yadong%Zhao%NULL%1,               Jiangtao%Liu%NULL%1,               Xiaotao%He%NULL%1,               Bo%Wang%NULL%0,               Shihua%Fu%NULL%1,               Jun%Yan%NULL%1,               Jingping%Niu%NULL%1,               Ji%Zhou%NULL%1,               Bin%Luo%luob@lzu.edu.cn%1]'
$ws.Range('E3').Value = '[Marcos Felipe Falcão%Sobral%NULL%1,               Gisleia Benini%Duarte%NULL%1,               Ana Iza Gomes%da Penha Sobral%NULL%1,               Marcelo Luiz Monteiro%Marinho%NULL%1,               André%de Souza Melo%NULL%1]'
$ws.Range('E4').Value = '[Dai%Su%NULL%1,               Yingchun%Chen%NULL%1,               Kevin%He%NULL%1,               Tao%Zhang%NULL%0,               Min%Tan%NULL%1,               Yunfan%Zhang%NULL%1,               Xingyu%Zhang%NULL%1,                D.%Su%null%1,                Y.% Chen%null%1,                K.% He%null%1,                T.% Zhang%null%1,                M.% Tan%null%1,                Y.% Zhang%null%1,                X. % Zhang%null%1,             D.%Su%null%1,             Y.% Chen%null%1,             K.% He%null%1,             T.% Zhang%null%1,             M.% Tan%null%1,             Y.% Zhang%null%1,             X. % Zhang%null%1]'
$ws.Range('E5').Value = '[Yu%Wu%NULL%1,               Wenzhan%Jing%NULL%1,               Jue%Liu%NULL%1,               Qiuyue%Ma%NULL%1,               Jie%Yuan%NULL%1,               Yaping%Wang%NULL%1,               Min%Du%NULL%1,               Min%Liu%NULL%0]'
$ws.Range('E6').Value = '[Yasir%Rehman%NULL%1,               Nadia%Rehman%NULL%1]'
$ws.Range('E7').Value = '[Cui%Guo%NULL%1,               Yacong%Bo%NULL%1,               Changqing%Lin%NULL%1,               Hao Bi%Li%NULL%1,               Yiqian%Zeng%NULL%1,               Yumiao%Zhang%NULL%1,               Md Shakhaoat%Hossain%NULL%1,               Jimmy W.M.%Chan%NULL%1,               David W.%Yeung%NULL%1,               Kin-on%Kwok%NULL%1,               Samuel Y.S.%Wong%NULL%1,               Alexis K.H.%Lau%NULL%1,               Xiang Qian%Lao%NULL%1]'
$ws.Range('E8').Value = '[Abu Reza Md. Towfiqul%Islam%towfiq_dm@brur.ac.bd%1,               Md.%Hasanuzzaman%NULL%2,               Md.%Hasanuzzaman%NULL%0,               Mashura%Shammi%NULL%1,               Roquia%Salam%NULL%1,               Md.%Bodrud-Doza%NULL%1,               Md. Mostafizur%Rahman%rahmanmm@juniv.edu%1,               Md. Abdul%Mannan%NULL%1,               Saleemul%Huq%NULL%1]'
$ws.Range('E9').Value = '[Ying%Jiang%yjiang1@vt.edu%1,               Jing%Xu%xujingjingalpha@sohu.com%0]'
$ws.Range('E10').Value = '[Yeran%Sun%NULL%1,               Xuke%Hu%NULL%1,               Jing%Xie%NULL%0]'
$ws.Range('E11').Value = '[Ploutarchos%Tzampoglou%NULL%1,               Dimitrios%Loukidis%NULL%1]'
$ws.Range('E12').Value = '[Daniel%Fernández%NULL%1,               Iago%Giné-Vázquez%NULL%1,               Ivy%Liu%NULL%1,               Recai%Yucel%NULL%1,               Marta%Nai Ruscone%NULL%1,               Marianthi%Morena%NULL%1,               Víctor Gerardo%García%NULL%1,               Josep Maria%Haro%NULL%1,               William%Pan%NULL%1,               Stefanos%Tyrovolas%NULL%1]'
